$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("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", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found: $found"
$rng.Collapse(0)
$rng.Text = "`r`rNew Classic Token Finance Mlops Project:`r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"
Write-Output $d.Content.Text.Length
